$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 is a text code ("004" -> "001"); force text so leading zeros are kept,
# then clear the temporary number-format style so no stray style is left behind.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

# Dates stored as plain text strings
$ws.Range("M2").Value = "2020-12-16 00:00:00"
$ws.Range("N2").Value = "2019-12-31 00:00:00"

# Numeric figures
$ws.Range("O2").Value = 1967653210.37
$ws.Range("P2").Value = 301931332.31
$ws.Range("Q2").Value = 215506807.88
$ws.Range("R2").Value = 79.23902244990001
$ws.Range("S2").Value = 151775981.35
$ws.Range("T2").Value = -28.5086973499
$ws.Range("U2").Value = 358016503.26
$ws.Range("V2").Value = -4.6636011729
$ws.Range("W2").Value = 882248564.09
$ws.Range("X2").Value = 257258612.44
$ws.Range("Y2").Value = 11.5134678086
$ws.Range("Z2").Value = 32240475.79
$ws.Range("AA2").Value = 1040.031745243
$ws.Range("AB2").Value = 1085404646.28
$ws.Range("AC2").Value = 19.8729360282
$ws.Range("AD2").Value = 15.4201717342
$ws.Range("AE2").Value = 10.3760658959
$ws.Range("AF2").Value = 168.0203343057
$ws.Range("AG2").Value = 44.8376044844
